$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column H:
#  - H1 holds a date (serial 42402 = 2016-02-02) formatted as a date (built-in
#    numFmtId 14, "mm-dd-yy").
#  - H2 is blank but carries the same date style.
#  - H3, H6, H7 hold a single-space text value and carry the same date style.
#  - H4, H5 hold the same single-space text value but with the default style.

$ws.Range("H1").Value = 42402
$ws.Range("H1").NumberFormat = "mm-dd-yy"

# Re-use H1's style (instead of re-creating an identical one) by copying its
# formatting onto the other date-styled cells.
$ws.Range("H1").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H7").PasteSpecial(-4122)

$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

$ws.Range("H7").Select()
